$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 21998
$ws.Range("E2").Value = 2489
$ws.Range("F2").Value = 2790
$ws.Range("G2").Value = 2289
$ws.Range("H2").Value = 1743
$ws.Range("I2").Value = 1635
$ws.Range("J2").Value = 108
$ws.Range("K2").Value = 29325
$ws.Range("L2").Value = 15300
$ws.Range("M2").Value = 14025
$ws.Range("N2").Value = 12471
$ws.Range("O2").Value = 1554
$ws.Range("P2").Value = 299
$ws.Range("Q2").Value = 2157
$ws.Range("R2").Value = -1206
$ws.Range("S2").Value = -1120
$ws.Range("T2").Value = 1550
$ws.Range("U2").Value = 607
$ws.Range("V2").Value = 10311
$ws.Range("W2").Value = 11.31
$ws.Range("X2").Value = 7.92
$ws.Range("Y2").Value = 13.8
$ws.Range("Z2").Value = 5.9
$ws.Range("AA2").Value = 109.09
$ws.Range("AB2").Value = 4240.78
$ws.Range("AC2").Value = 2736
$ws.Range("AD2").Value = 23.66
$ws.Range("AE2").Value = 23684
$ws.Range("AF2").Value = 2.73
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 0.93
$ws.Range("AI2").Value = 19.33
$ws.Range("AJ2").Value = 59749690

# Row 3
$ws.Range("D3").Value = 23824
$ws.Range("E3").Value = 2993
$ws.Range("F3").Value = 2993
$ws.Range("G3").Value = 2518
$ws.Range("H3").Value = 1771
$ws.Range("I3").Value = 1623
$ws.Range("J3").Value = 148
$ws.Range("K3").Value = 29988
$ws.Range("L3").Value = 14727
$ws.Range("M3").Value = 15262
$ws.Range("N3").Value = 13787
$ws.Range("O3").Value = 1475
$ws.Range("P3").Value = 299
$ws.Range("Q3").Value = 4821
$ws.Range("R3").Value = -3022
$ws.Range("S3").Value = -1439
$ws.Range("T3").Value = 1350
$ws.Range("U3").Value = 3471
$ws.Range("V3").Value = 10070
$ws.Range("W3").Value = 12.56
$ws.Range("X3").Value = 7.43
$ws.Range("Y3").Value = 12.36
$ws.Range("Z3").Value = 5.97
$ws.Range("AA3").Value = 96.48999999999999
$ws.Range("AB3").Value = 4674.35
$ws.Range("AC3").Value = 2716
$ws.Range("AD3").Value = 27.38
$ws.Range("AE3").Value = 26255
$ws.Range("AF3").Value = 2.83
$ws.Range("AG3").Value = 600
$ws.Range("AH3").Value = 0.8100000000000001
$ws.Range("AI3").Value = 19.41
$ws.Range("AJ3").Value = 59760910

# Row 4
$ws.Range("D4").Value = 1280
$ws.Range("E4").Value = -18
$ws.Range("F4").Value = 3262
$ws.Range("G4").Value = -15
$ws.Range("H4").Value = 2490
$ws.Range("I4").Value = 2367
$ws.Range("J4").Value = 124
$ws.Range("K4").Value = 29089
$ws.Range("L4").Value = 12040
$ws.Range("M4").Value = 17049
$ws.Range("N4").Value = 15995
$ws.Range("O4").Value = 1054
$ws.Range("P4").Value = 300
$ws.Range("Q4").Value = 3406
$ws.Range("R4").Value = -430
$ws.Range("S4").Value = -3089
$ws.Range("T4").Value = 2371
$ws.Range("U4").Value = 1035
$ws.Range("V4").Value = 7469
$ws.Range("W4").Value = -1.4
$ws.Range("X4").Value = 194.62
$ws.Range("Y4").Value = 15.89
$ws.Range("Z4").Value = 8.43
$ws.Range("AA4").Value = 70.62
$ws.Range("AB4").Value = 5339.82
$ws.Range("AC4").Value = 3952
$ws.Range("AD4").Value = 10.57
$ws.Range("AE4").Value = 30274
$ws.Range("AF4").Value = 1.38
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 1.44
$ws.Range("AI4").Value = 13.39
$ws.Range("AJ4").Value = 60084350

# Row 5
$ws.Range("D5").Value = 1127
$ws.Range("E5").Value = 49
$ws.Range("F5").Value = 49
$ws.Range("G5").Value = -400
$ws.Range("H5").Value = 14864
$ws.Range("I5").Value = 14803
$ws.Range("J5").Value = 61
$ws.Range("K5").Value = 43028
$ws.Range("L5").Value = 11849
$ws.Range("M5").Value = 31179
$ws.Range("N5").Value = 18479
$ws.Range("O5").Value = 12700
$ws.Range("P5").Value = 313
$ws.Range("Q5").Value = -88
$ws.Range("R5").Value = 591
$ws.Range("S5").Value = 144
$ws.Range("T5").Value = 581
$ws.Range("U5").Value = -669
$ws.Range("V5").Value = 5997
$ws.Range("W5").Value = 4.35
$ws.Range("X5").Value = 1318.34
$ws.Range("Y5").Value = 85.88
$ws.Range("Z5").Value = 41.22
$ws.Range("AA5").Value = 38.01
$ws.Range("AB5").Value = 14492.55
$ws.Range("AC5").Value = 34974
$ws.Range("AD5").Value = 0.76
$ws.Range("AE5").Value = 30718
$ws.Range("AF5").Value = 0.87
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 2.25
$ws.Range("AI5").Value = 1.39
$ws.Range("AJ5").Value = 62644053

# Row 6
$ws.Range("D6").Value = 19977
$ws.Range("E6").Value = 2305
$ws.Range("F6").Value = 2305
$ws.Range("G6").Value = 2409
$ws.Range("H6").Value = 1183
$ws.Range("I6").Value = 563
$ws.Range("K6").Value = 43421
$ws.Range("L6").Value = 11655
$ws.Range("M6").Value = 31766
$ws.Range("N6").Value = 18761
$ws.Range("P6").Value = 313
$ws.Range("Q6").Value = 3505
$ws.Range("R6").Value = -1224
$ws.Range("S6").Value = -1822
$ws.Range("T6").Value = 1671
$ws.Range("U6").Value = 1834
$ws.Range("V6").Value = 4650
$ws.Range("W6").Value = 11.54
$ws.Range("X6").Value = 5.92
$ws.Range("Y6").Value = 3.02
$ws.Range("Z6").Value = 2.74
$ws.Range("AA6").Value = 36.69
$ws.Range("AB6").Value = 5815.83
$ws.Range("AC6").Value = 899
$ws.Range("AD6").Value = 19.02
$ws.Range("AE6").Value = 31187
$ws.Range("AF6").Value = 0.55
$ws.Range("AG6").Value = 650
$ws.Range("AH6").Value = 3.8
$ws.Range("AI6").Value = 38.2
$ws.Range("AJ6").Value = 62644293

# Clear estimate rows 7-9 (columns D through AJ), keep A/B/C
$ws.Range("D7:AJ9").ClearContents()
